$wb = $excel.ActiveWorkbook

# The new handoff was generated for e2e\b.md: it is now "Ready for handoff"
# with a fresh handoff package (b.*.xlf) and timestamp, and a note that the
# previously handed-back file is stale relative to the newest source.

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7d0fd6d22951f2572bd4cfcf520aec50f607d8d1/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ab667e5a8d5c99f8ca6ce88c82aaedc5ee8373e8/e2e/b.md."

# ---- Overview sheet: update the b.md row (row 3) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-04 20:43:47"

# ---- zh-cn sheet: update the b.md row (row 3) ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-04 20:43:43"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39 + 1/6

# ---- de-de sheet: update the b.md row (row 3) ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-04 20:43:47"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39 + 1/6
